# Enter grading points in column E (Total Points) for the Customer Class
# and Product Class sections, mirroring the full marks already recorded
# in column D (Points for grading).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Customer Class section (rows 3-6)
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Product Class section (rows 10-14)
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Move selection to E15 (grader's next stop)
$ws.Range("E15").Select()
